# CreateResultsDF change: the "NA" filler columns (K:AC) that used to sit on
# the same rows as the real equipment/system data (rows 2-6) are relocated to
# their own dedicated block of rows beneath all the data rows, preventing a
# column mismatch between the two concatenated frames.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the K2:AC6 "NA" block down to start at K12 (rows 2-6 -> rows 12-16),
# clearing it out of rows 2-6 in the process.
$src = $ws.Range("K2:AC6")
$dst = $ws.Range("K12")
$src.Cut($dst)

# Append a further 5 rows (17-21) of the same "NA" filler block.
$ws.Range("K17:AC21").Value = "NA"
